# 24 apr cdc update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Revise existing daily case counts (CDC revisions to prior days' totals)
$ws.Range("B43").Value = 80
$ws.Range("B46").Value = 214
$ws.Range("B47").Value = 279
$ws.Range("B48").Value = 423
$ws.Range("B55").Value = 3487
$ws.Range("B57").Value = 7038
$ws.Range("B62").Value = 33404
$ws.Range("B63").Value = 44183
$ws.Range("B74").Value = 277205
$ws.Range("B75").Value = 304826
$ws.Range("B78").Value = 395011

# Append new daily rows for 21-23 Apr 2020
$ws.Range("A92").Value = 43942
$ws.Range("B92").Value = 802583

$ws.Range("A93").Value = 43943
$ws.Range("B93").Value = 828441

$ws.Range("A94").Value = 43944
$ws.Range("B94").Value = 865585

# Match date formatting used by the rest of column A
$ws.Range("A92:A94").NumberFormat = "[$-409]d\-mmm\-yyyy;@"
